# Update "想去人数" (number of people wanting to go) figures on both the
# "展览" and "全部类型" sheets, which hold duplicate copies of the same data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 8613
    $ws.Range("F3").Value = 194
    $ws.Range("F4").Value = 395
}
